# Update cryptos list price (D) and 1h volume-change (E) columns.
# A leading apostrophe forces Excel to store the price as text
# (matching the source data, which stores these as literal strings,
# e.g. "1.00" / "3.284.26" rather than as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.460.23"
$ws.Range("E2").Value = "  +1.79%  "

$ws.Range("D3").Value = "'3.284.36"
$ws.Range("E3").Value = "  +5.23%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'249.19"
$ws.Range("E5").Value = "  +4.57%  "

$ws.Range("D6").Value = "'617.86"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  -0.93%  "

$ws.Range("D8").Value = "'0.382"
$ws.Range("E8").Value = "  -2.69%  "

$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").Value = "'3.284.26"
$ws.Range("E10").Value = "  +5.30%  "

$ws.Range("D11").Value = "'0.781"
$ws.Range("E11").Value = "  -6.76%  "

$ws.Range("D12").Value = "'0.198"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").Value = "'95.987.41"
$ws.Range("E13").Value = "  +1.71%  "

$ws.Range("D14").Value = "'0.0000244"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").Value = "'34.99"
$ws.Range("E15").Value = "  +0.64%  "

$ws.Range("D16").Value = "'3.875.04"
$ws.Range("E16").Value = "  +4.75%  "

$ws.Range("D17").Value = "'5.49"
$ws.Range("E17").Value = "  +3.94%  "

$ws.Range("D18").Value = "'3.281.08"
$ws.Range("E18").Value = "  +5.06%  "

$ws.Range("D19").Value = "'3.57"
$ws.Range("E19").Value = "  -3.05%  "

$ws.Range("D20").Value = "'14.87"
$ws.Range("E20").Value = "  -0.95%  "

$ws.Range("D21").Value = "'480.03"
$ws.Range("E21").Value = "  +6.50%  "

$ws.Range("D22").Value = "'5.79"
$ws.Range("E22").Value = "  -2.31%  "

$ws.Range("D23").Value = "'0.0000204"
$ws.Range("E23").Value = "  +2.31%  "

$ws.Range("D24").Value = "'9.17"
$ws.Range("E24").Value = "  +1.79%  "

$ws.Range("D25").Value = "'5.59"
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("D26").Value = "'87.41"
$ws.Range("E26").Value = "  +1.65%  "

$ws.Range("D27").Value = "'12.00"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").Value = "'3.458.64"
$ws.Range("E28").Value = "  +5.17%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.38%  "

$ws.Range("D30").Value = "'0.181"
$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("D31").Value = "'0.238"
$ws.Range("E31").Value = "  -8.23%  "

$ws.Range("E32").Value = "  +1.00%  "

$ws.Range("E33").Value = "  -2.69%  "

$ws.Range("D34").Value = "'9.19"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").Value = "'27.11"
$ws.Range("E35").Value = "  +4.19%  "

$ws.Range("D36").Value = "'7.36"
$ws.Range("E36").Value = "  -7.49%  "

$ws.Range("D37").Value = "'0.150"
$ws.Range("E37").Value = "  -6.39%  "

$ws.Range("D38").Value = "'1.92"
$ws.Range("E38").Value = "  +0.54%  "

$ws.Range("D39").Value = "'24.68"
$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").Value = "'493.14"
$ws.Range("E40").Value = "  +3.11%  "

$ws.Range("D41").Value = "'0.445"
$ws.Range("E41").Value = "  -2.93%  "

$ws.Range("D42").Value = "'1.27"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("D45").Value = "'3.45"
$ws.Range("E45").Value = "  -4.98%  "

$ws.Range("D46").Value = "'0.787"
$ws.Range("E46").Value = "  +13.26%  "

$ws.Range("D47").Value = "'161.14"
$ws.Range("E47").Value = "  -0.28%  "

$ws.Range("D48").Value = "'1.89"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("D49").Value = "'45.15"
$ws.Range("E49").Value = "  +2.77%  "

$ws.Range("D50").Value = "'1.35"
$ws.Range("E50").Value = "  +2.81%  "

$ws.Range("D51").Value = "'4.49"
$ws.Range("E51").Value = "  +1.49%  "

# Rows 43/44: coin identity swap (USDe <-> dogwifhat), each also getting
# refreshed price/volume-change values.
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'3.24"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  -0.01%  "
